$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.51"
$ws.Range("E2").Value = "'6.24%"
$ws.Range("G2").Value = "'15"
$ws.Range("D3").Value = "'27.11"
$ws.Range("E3").Value = "'1.32%"
$ws.Range("G3").Value = "'15"
$ws.Range("D4").Value = "'4.908"
$ws.Range("E4").Value = "'4.90%"
$ws.Range("G4").Value = "'15"
$ws.Range("D5").Value = "'0.06358"
$ws.Range("E5").Value = "'4.06%"
$ws.Range("G5").Value = "'15"
$ws.Range("D6").Value = "'6.947"
$ws.Range("E6").Value = "'3.60%"
$ws.Range("G6").Value = "'15"
$ws.Range("E7").Value = "'6.00%"
$ws.Range("G7").Value = "'15"
$ws.Range("D8").Value = "'0.8852"
$ws.Range("E8").Value = "'4.07%"
$ws.Range("G8").Value = "'15"
$ws.Range("D9").Value = "'0.9485"
$ws.Range("E9").Value = "'4.33%"
$ws.Range("G9").Value = "'15"
$ws.Range("D10").Value = "'0.1475"
$ws.Range("E10").Value = "'4.61%"
$ws.Range("G10").Value = "'15"
$ws.Range("D11").Value = "'0.05135"
$ws.Range("E11").Value = "'3.20%"
$ws.Range("G11").Value = "'15"
$ws.Range("D12").Value = "'0.07475"
$ws.Range("E12").Value = "'5.38%"
$ws.Range("G12").Value = "'15"
$ws.Range("E13").Value = "'0.49%"
$ws.Range("G13").Value = "'15"
$ws.Range("D14").Value = "'0.09047"
$ws.Range("E14").Value = "'-0.07%"
$ws.Range("G14").Value = "'15"
$ws.Range("D15").Value = "'0.001557"
$ws.Range("E15").Value = "'1.78%"
$ws.Range("G15").Value = "'15"
$ws.Range("D16").Value = "'0.0006260"
$ws.Range("E16").Value = "'0.74%"
$ws.Range("G16").Value = "'15"
$ws.Range("D17").Value = "'0.005782"
$ws.Range("E17").Value = "'-3.02%"
$ws.Range("G17").Value = "'15"
$ws.Range("D18").Value = "'3.484"
$ws.Range("E18").Value = "'1.01%"
$ws.Range("G18").Value = "'15"
$ws.Range("D19").Value = "'2.296"
$ws.Range("E19").Value = "'5.91%"
$ws.Range("G19").Value = "'15"
$ws.Range("E20").Value = "'0.86%"
$ws.Range("G20").Value = "'15"
$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").Value = "'1.53%"
$ws.Range("G21").Value = "'15"
$ws.Range("D22").Value = "'3.850"
$ws.Range("E22").Value = "'-5.64%"
$ws.Range("G22").Value = "'15"
$ws.Range("D23").Value = "'0.04325"
$ws.Range("E23").Value = "'2.14%"
$ws.Range("G23").Value = "'15"
$ws.Range("E24").Value = "'-0.14%"
$ws.Range("G24").Value = "'15"
$ws.Range("D25").Value = "'0.003634"
$ws.Range("E25").Value = "'-10.51%"
$ws.Range("G25").Value = "'15"
$ws.Range("D26").Value = "'0.0001198"
$ws.Range("E26").Value = "'-0.18%"
$ws.Range("G26").Value = "'15"
$ws.Range("D27").Value = "'0.0001691"
$ws.Range("E27").Value = "'-12.71%"
$ws.Range("G27").Value = "'15"
$ws.Range("G28").Value = "'15"
$ws.Range("G29").Value = "'15"
$ws.Range("G30").Value = "'15"
$ws.Range("G31").Value = "'15"
$ws.Range("G32").Value = "'15"
$ws.Range("G33").Value = "'15"
$ws.Range("G34").Value = "'15"
$ws.Range("G35").Value = "'15"
$ws.Range("G36").Value = "'15"
$ws.Range("G37").Value = "'15"
$ws.Range("G38").Value = "'15"
$ws.Range("G39").Value = "'15"
$ws.Range("D40").Value = "'0.04062"
$ws.Range("E40").Value = "'3.01%"
$ws.Range("G40").Value = "'15"
$ws.Range("D41").Value = "'0.006624"
$ws.Range("E41").Value = "'59.32%"
$ws.Range("G41").Value = "'15"
$ws.Range("D42").Value = "'0.1165"
$ws.Range("E42").Value = "'4.60%"
$ws.Range("G42").Value = "'15"
$ws.Range("D43").Value = "'0.002346"
$ws.Range("E43").Value = "'11.20%"
$ws.Range("G43").Value = "'15"
$ws.Range("D44").Value = "'0.01251"
$ws.Range("E44").Value = "'8.66%"
$ws.Range("G44").Value = "'15"
$ws.Range("D45").Value = "'0.00005239"
$ws.Range("E45").Value = "'3.30%"
$ws.Range("G45").Value = "'15"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("G46").Value = "'15"
$ws.Range("D47").Value = "'2.248"
$ws.Range("E47").Value = "'795.44%"
$ws.Range("G47").Value = "'15"
$ws.Range("D48").Value = "'0.02257"
$ws.Range("E48").Value = "'6.37%"
$ws.Range("G48").Value = "'15"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("G49").Value = "'15"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("G50").Value = "'15"
$ws.Range("G51").Value = "'15"
